# Trace report refresh: new search results (06/15/2023), replacing the
# previous 05/12/2023 pull. Row 3 ("1 CO" banner) is newly inserted above
# the header row, and the data block grows from 6 to 8 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Make room: insert a new row above the header (old row 3) so the
#    header + existing data shift down by one. The header row keeps its
#    row-level bold/customFormat styling automatically.
# ---------------------------------------------------------------------
$ws.Rows("3:3").Insert()

# ---------------------------------------------------------------------
# 1b. Drop the AutoFilter and register the worksheet-level sortState that
#     records how the trace data is ordered (custom lists: state "CO"
#     first, then event "Placed Actual" before "Placed Construct"). This
#     must run before the data rows are overwritten below, since Apply()
#     performs a real sort of whatever is currently in the range.
# ---------------------------------------------------------------------
$ws.AutoFilterMode = $false

$sf = $ws.Sort
$sf.SortFields.Clear()
$f1 = $sf.SortFields.Add($ws.Range("D5:D12"))
$f1.CustomOrder = "CO"
$f2 = $sf.SortFields.Add($ws.Range("H5:H12"))
$f2.CustomOrder = "Placed Actual,Placed Construct"
$sf.SetRange($ws.Range("A4:O12"))
$sf.Header = 1
$sf.Apply()

# ---------------------------------------------------------------------
# 2. Title / summary rows
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Description unknown, completed 06/15/2023 06:00:44 EDT, by WPJTOWN1.The search returned: 8 events."

$ws.Range("A2").Value = "5 On Hand"
$ws.Range("A2").Interior.Color = 16776960

$ws.Range("A3").Value = "1 CO"
$ws.Range("A3").Interior.Color = 5874847

# ---------------------------------------------------------------------
# 3. Header row (row 4) - values only, style already carried over by the
#    row insert above.
# ---------------------------------------------------------------------
$header = New-Object 'object[,]' 1,15
$header[0,0]  = "Initial"
$header[0,1]  = "Number"
$header[0,2]  = "Location City"
$header[0,3]  = "State"
$header[0,4]  = "Month"
$header[0,5]  = "Day"
$header[0,6]  = "Time"
$header[0,7]  = "Event"
$header[0,8]  = "Train ID"
$header[0,9]  = "Destination City"
$header[0,10] = "State"
$header[0,11] = "Gross Weight"
$header[0,12] = "Tare Weight"
$header[0,13] = "Net Weight"
$header[0,14] = "Car_no"
$ws.Range("A4:O4").Value = $header

# ---------------------------------------------------------------------
# 4. Data rows 5-12
# ---------------------------------------------------------------------
$data = New-Object 'object[,]' 8,15
# Row 5 - ITFX 9725
$data[0,0]="ITFX";  $data[0,1]=9725;   $data[0,2]="JOHNSTOWN"; $data[0,3]="CO"; $data[0,4]=6; $data[0,5]=1;  $data[0,6]=1812; $data[0,7]="Placed Actual";     $data[0,8]=""; $data[0,9]="JOHNSTOWN"; $data[0,10]="CO"; $data[0,11]=202700; $data[0,12]=0;     $data[0,13]=202700; $data[0,14]="ITFX9725"
# Row 6 - ITFX 9728
$data[1,0]="ITFX";  $data[1,1]=9728;   $data[1,2]="JOHNSTOWN"; $data[1,3]="CO"; $data[1,4]=6; $data[1,5]=1;  $data[1,6]=1812; $data[1,7]="Placed Actual";     $data[1,8]=""; $data[1,9]="JOHNSTOWN"; $data[1,10]="CO"; $data[1,11]=202950; $data[1,12]=0;     $data[1,13]=202950; $data[1,14]="ITFX9728"
# Row 7 - MWCX 102555
$data[2,0]="MWCX";  $data[2,1]=102555; $data[2,2]="JOHNSTOWN"; $data[2,3]="CO"; $data[2,4]=6; $data[2,5]=1;  $data[2,6]=1811; $data[2,7]="Placed Actual";     $data[2,8]=""; $data[2,9]="LOVELAND";  $data[2,10]="CO"; $data[2,11]=281100; $data[2,12]=73600; $data[2,13]=207500; $data[2,14]="MWCX102555"
# Row 8 - MWCX 102276
$data[3,0]="MWCX";  $data[3,1]=102276; $data[3,2]="JOHNSTOWN"; $data[3,3]="CO"; $data[3,4]=6; $data[3,5]=12; $data[3,6]=1304; $data[3,7]="Placed Actual";     $data[3,8]=""; $data[3,9]="LOVELAND";  $data[3,10]="CO"; $data[3,11]=280350; $data[3,12]=78900; $data[3,13]=201450; $data[3,14]="MWCX102276"
# Row 9 - MWCX 102166
$data[4,0]="MWCX";  $data[4,1]=102166; $data[4,2]="JOHNSTOWN"; $data[4,3]="CO"; $data[4,4]=6; $data[4,5]=12; $data[4,6]=1304; $data[4,7]="Placed Actual";     $data[4,8]=""; $data[4,9]="LOVELAND";  $data[4,10]="CO"; $data[4,11]=282400; $data[4,12]=82000; $data[4,13]=200400; $data[4,14]="MWCX102166"
# Row 10 - MWCX 102330
$data[5,0]="MWCX";  $data[5,1]=102330; $data[5,2]="LOVELAND";  $data[5,3]="CO"; $data[5,4]=6; $data[5,5]=12; $data[5,6]=1045; $data[5,7]="Junction Received"; $data[5,8]="BNSF"; $data[5,9]="LOVELAND";  $data[5,10]="CO"; $data[5,11]=284850; $data[5,12]=79300; $data[5,13]=205550; $data[5,14]="MWCX102330"
# Row 11 - MWCX 102328
$data[6,0]="MWCX";  $data[6,1]=102328; $data[6,2]="VALDOSTA";  $data[6,3]="GA"; $data[6,4]=6; $data[6,5]=14; $data[6,6]=2248; $data[6,7]="Arrive In-Transit"; $data[6,8]=""; $data[6,9]="LOVELAND";  $data[6,10]="CO"; $data[6,11]=280550; $data[6,12]=79500; $data[6,13]=201050; $data[6,14]="MWCX102328"
# Row 12 - MWCX 102553
$data[7,0]="MWCX";  $data[7,1]=102553; $data[7,2]="MEMPHIS";   $data[7,3]="TN"; $data[7,4]=6; $data[7,5]=14; $data[7,6]=1300; $data[7,7]="Junction Received"; $data[7,8]="NS"; $data[7,9]="LOVELAND";  $data[7,10]="CO"; $data[7,11]=281050; $data[7,12]=73400; $data[7,13]=207650; $data[7,14]="MWCX102553"

$ws.Range("A5:O12").Value = $data

# ---------------------------------------------------------------------
# 5. Row-level fill colors mirroring the event highlighting
#    (A:N per row; O/"Car_no" column is left unfilled, as in the source)
# ---------------------------------------------------------------------
$ws.Range("A5:N9").Interior.Color = 16776960    # yellow - Placed Actual rows
$ws.Range("A10:N10").Interior.Color = 5874847   # olive - Junction Received (CO)
$ws.Range("A11:N11").Interior.Color = 255       # red - Arrive In-Transit (GA)
# Row 12 keeps the default (no fill)

# ---------------------------------------------------------------------
# 6. Defined name range for the (now-inactive) filter database
# ---------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Test_format_trace!`$A`$4:`$O`$12"
    }
}

# ---------------------------------------------------------------------
# 8. Selection / view state
# ---------------------------------------------------------------------
$ws.Range("K5:K12").Select()
